$d = $word.ActiveDocument

# --- 1. Update the title paragraph (date + paper title), keeps the <w:br/> split ---
$d.Content.Find.Execute("המאמר היומי של מייק - 27.12.24:", $true, $false, $false, $false, $false,
                         $true, 1, $false, "המאמר היומי של מייק - 26.12.24:", 2)
$d.Content.Find.Execute("Position: Future Directions in the Theory of Graph Machine Learning", $true, $false, $false, $false, $false,
                         $true, 1, $false, "RL for Consistency Models: Faster Reward Guided Text-to-Image Generation", 2)

# --- 2. Replace the body text of paragraphs 2-13 (1-indexed) with the new review content ---
function Set-ParaText($index, $text) {
    $p = $d.Paragraphs.Item($index)
    $r = $p.Range
    $r.MoveEnd(1, -1)
    $r.Text = $text
}

Set-ParaText 2 ' מזמן לא סקרתי מאמרים על מודלי דיפוזיה אז אחרי שנתקלתי במאמר הנחמד המשלב מודלי דיפוזיה גנרטיביים עם למידה עם חיזוקים (Reinforcement Learning או RL בקצרה), לא היו לי ספקות שזה הולך להיות המאמר המסוקר. כאמור המאמר פיתח שיטת אימון מודל של דיפוזיה גנרטיבי מסוג Consistency Model או CM.'

Set-ParaText 3 'קודם כל נשאלת השאלה למה צריך לאמן מודלי דיפוזיה גנרטיביים עם שיטות הלקוחות מעולם RL. הרי יש לנו שיטות סטנדרטיות יותר לאימון של מודלי דיפוזיה שהצליחו להביא לנו מודלים בעלי ביצועים מרשימים (בגנרוט תמונות מטקסט). אתם בטח יודעים שאימון מודלי דיפוזיה לגנרוט תמונות זה דבר לא זול ודורש לא מעט זמן ושימוש RL לאימון (או fine-tune) של מודלי דיפוזיה יכול לחסוך לנו זמן במקרים שאנו צריכים לאמן מודל דיפוזיה ייעודי (למשל לדומיין נישתי) '

Set-ParaText 4 'אחת הדוגמאות למשימה כזו היא אימון מודל ליצירת תמונות מפרומפט (תיאור טקסטואלי) כאשר יש בידינו פונקציה המשערכת את התאמת התמונה לפרומפט. אתם כבר יכולים לנחש שפונקציה זו תשרת לנו בתור פונקצית תגמול (reward function). '

Set-ParaText 5 'כבר הזכרתי שהמאמר משלב שיטה חדשה (יחסית) לאימון מודלי דיפוזיה הנקראת CM ושיטה זו (שהומצאה על ידי איליה סלוצקב ושות'') מאפשרת גנרוט יותר מהיר של מודלי דיפוזיה גנרטיביים. בגדול מאוד שיטה זו מנסה לאמן מודל שאוכף עקביות בין התמונות המשוחזרות על ידי המודל מתמונות מורעשות עם עוצמות שונות רעש. כלומר לוקחים תמונה, מרעישים אותה עם רעש (בד״כ גאוסי) עם שונויות שונות ומאמנים מודל להחזיר את אותה התמונה הנקייה (עקביות לשמה). '

Set-ParaText 6 'למה השיטה הזו מאפשרת גנרוט יותר מהיר של תמונות? כי בגדול היא מאפשרת לגנרט תמונה נקייה מרעש באיטרציה אחת בלבד (ככה המודל מאומן). במציאות עושים את זה בכמה איטרציות (מספר קטן). מתחילים מרעש, מגנרטים את התמונה ממנו, מוסיפים פחות רעש לתמונה המגונרטת, מגנרטים מהתמונה המורעשת שוב וממשיכים ככה כמה איטרציות (עשרות בודדת). זה מאפשר לזרז את תהליך הגנרוט כי מודלי דיפוזיה סטנדרטיים צריכים מאות איטרציות בד״כ.'

Set-ParaText 7 'אוקיי, אחרי הקדמה ארוכה נעבור לתיאור של מה שעשו במאמר. המחברים הגדירו Markov Decision Process c או MDP המתאר תהליך גנרוט של תמונה (או כל דאטה אחר למעשה). כאמור פונקציה תגמול ניתנת לנו והיא מודדת מידת התאמה של התמונה המגונרטת לפרומפט. המאמר מגדיר:'

Set-ParaText 8 'המצב s_t בתור שלישיה התמונה מגונרטת באיטרציה t, עוצמת הרעש והפרומפט c'

Set-ParaText 9 'הפעולה a_t היא התמונה באיטרציה t + 1'

Set-ParaText 10 'הפוליסי היא זו פונקצית התפלגות מותנית של תמונה מאיטרציה t+1 בהינתן התמונה המגונרטת מאיטרציה t בתוספת רעש'

Set-ParaText 11 'המצב המתחלתי הוא רעש גאוסי סטנדרטי ופונקציית תגמול נתונה לנו'

Set-ParaText 12 'אחרי שהגדרנו את ה-MDP של תהליך גנרוט התמונה אנו יכולים להשתמש בשיטה DPO או Direct Preference Optimization לאימון פונקצית עקביות (= המודל שאנו מאמנים). למעשה DPO מאמן מודל הממקסם את פונקצית התגמול תוך כדי הגבלת של גודל עדכון פרמטרי המודל בכל איטרציה (הומצא על ג''ו שולמן ה-CTO של OpenAI לשעבר).'

Set-ParaText 13 'המאמר גם טוען שאימון כזה הוא חסכוני מבחינת משאבי החישוב הנדרשים ויעיל מבחינה הדאטה (כלומר יכול לעבוד לדאטהסטים קטנים).'

# --- 3. Delete the old detailed GNN-review paragraphs (14 through 46, 1-indexed), ---
#        leaving only the final URL paragraph (now re-indexed to 14).
for ($i = 46; $i -ge 14; $i--) {
    $d.Paragraphs.Item($i).Range.Delete()
}

# --- 4. Update the arxiv link text in the final paragraph ---
$d.Content.Find.Execute("https://arxiv.org/abs/2402.02287", $true, $false, $false, $false, $false,
                         $true, 1, $false, "https://arxiv.org/abs/2404.03673", 2)
